$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.504.86'
$ws.Range('E2').Value = '  -0.49%  '

# Row 3
$ws.Range('D3').Value = '3.451.49'
$ws.Range('E3').Value = '  -0.05%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.98'
$ws.Range('E5').Value = '  -1.34%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.52'
$ws.Range('E6').Value = '  -1.20%  '

# Row 7
$ws.Range('E7').Value = '  +0.10%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  -1.61%  '

# Row 9
$ws.Range('D9').Value = '3.447.20'
$ws.Range('E9').Value = '  -0.12%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.132'
$ws.Range('E10').Value = '  -5.18%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.88'
$ws.Range('E11').Value = '  -0.87%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.418'
$ws.Range('E12').Value = '  -2.93%  '

# Row 13
$ws.Range('D13').Value = '4.051.97'
$ws.Range('E13').Value = '  +0.08%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.30'
$ws.Range('E14').Value = '  -4.28%  '

# Row 15
$ws.Range('E15').Value = '  -0.46%  '

# Row 16
$ws.Range('D16').Value = '66.721.92'
$ws.Range('E16').Value = '  -0.16%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000172'
$ws.Range('E17').Value = '  -2.35%  '

# Row 18
$ws.Range('D18').Value = '3.471.45'
$ws.Range('E18').Value = '  +0.53%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.97'
$ws.Range('E19').Value = '  -4.05%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('E20').Value = '  -2.03%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '378.95'
$ws.Range('E21').Value = '  -2.54%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.82'
$ws.Range('E22').Value = '  -0.91%  '

# Row 23
$ws.Range('E23').Value = '  +1.13%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.37'
$ws.Range('E24').Value = '  +0.69%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.997'
$ws.Range('E25').Value = '  -0.17%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.71'
$ws.Range('E26').Value = '  -0.70%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000120'
$ws.Range('E27').Value = '  -0.44%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.81'
$ws.Range('E28').Value = '  -4.52%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.176'
$ws.Range('E29').Value = '  +1.31%  '

# Row 30
$ws.Range('E30').Value = '  -0.05%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.36'
$ws.Range('E31').Value = '  +4.18%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.85'
$ws.Range('E32').Value = '  -4.65%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.00'
$ws.Range('E33').Value = '  -2.53%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.32'
$ws.Range('E34').Value = '  -5.62%  '

# Row 35
$ws.Range('E35').Value = '  +0.06%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.17'
$ws.Range('E36').Value = '  -1.56%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.58'
$ws.Range('E37').Value = '  +0.22%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.71'
$ws.Range('E38').Value = '  -1.68%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '29.43'
$ws.Range('E39').Value = '  +12.82%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.892'
$ws.Range('E40').Value = '  +2.05%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.65'
$ws.Range('E41').Value = '  -5.53%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.79'
$ws.Range('E42').Value = '  -4.21%  '

# Row 43
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.53'
$ws.Range('E43').Value = '  -2.42%  '

# Row 44
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.751.15'
$ws.Range('E44').Value = '  +0.79%  '

# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.40'
$ws.Range('E45').Value = '  -4.99%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0695'
$ws.Range('E46').Value = '  -3.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.65'
$ws.Range('E47').Value = '  -1.09%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.63'
$ws.Range('E48').Value = '  -5.34%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0294'
$ws.Range('E49').Value = '  -1.06%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '310.02'
$ws.Range('E50').Value = '  -4.67%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.830'
$ws.Range('E51').Value = '  -0.70%  '
